$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data shrank from 22 data rows to 20: drop the old trailing rows 21 and 22
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Delete()

# Column C got a bit wider (the engine's ColumnWidth setter adds a fixed
# ~0.8333 offset vs. the raw OOXML stored "width" unit, so back it out to
# land exactly on 27)
$ws.Columns.Item(3).ColumnWidth = 26.166666666666668

# Helper for the handful of values that look like dates (dd/mm/yyyy) and
# would otherwise get silently reinterpreted as date serials on assignment.
function Set-TextLiteral($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 3).Value = "Sabão Em Barra"
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = ""

# Row 3
$ws.Cells.Item(3, 1).Value = 20
$ws.Cells.Item(3, 2).Value = "PAPELARIA"
$ws.Cells.Item(3, 3).Value = "Borracha Branca"
$ws.Cells.Item(3, 4).Value = 1

# Row 4
$ws.Cells.Item(4, 3).Value = "Chave De Fenda"

# Row 5
$ws.Cells.Item(5, 1).Value = 28
$ws.Cells.Item(5, 2).Value = "MATERIALESCOLAR"
$ws.Cells.Item(5, 3).Value = "Caderno Espiral"
$ws.Cells.Item(5, 4).Value = 200

# Row 6
$ws.Cells.Item(6, 1).Value = 20
$ws.Cells.Item(6, 2).Value = "PAPELARIA"
$ws.Cells.Item(6, 3).Value = "Papel A4 Sulfite"
$ws.Cells.Item(6, 4).Value = 25

# Row 7
$ws.Cells.Item(7, 1).Value = 24
$ws.Cells.Item(7, 2).Value = "INFORMATICA"
$ws.Cells.Item(7, 3).Value = "Mouse Usb"

# Row 8
$ws.Cells.Item(8, 1).Value = 20
$ws.Cells.Item(8, 2).Value = "PAPELARIA"
$ws.Cells.Item(8, 3).Value = "Caneta Preta"
$ws.Cells.Item(8, 4).Value = 4
Set-TextLiteral 8 5 "10/12/2025"

# Row 9
$ws.Cells.Item(9, 3).Value = "Caneta Azul"
$ws.Cells.Item(9, 4).Value = 4
$ws.Cells.Item(9, 5).Value = ""

# Row 10
$ws.Cells.Item(10, 1).Value = 26
$ws.Cells.Item(10, 2).Value = "LIMPEZA"
$ws.Cells.Item(10, 3).Value = "Sabão Em Pó"
$ws.Cells.Item(10, 4).Value = 2

# Row 11
$ws.Cells.Item(11, 1).Value = 20
$ws.Cells.Item(11, 2).Value = "PAPELARIA"
$ws.Cells.Item(11, 3).Value = "Lápis"
$ws.Cells.Item(11, 4).Value = 4

# Row 12
$ws.Cells.Item(12, 1).Value = 18
$ws.Cells.Item(12, 2).Value = "ELETRICOS"
$ws.Cells.Item(12, 3).Value = "Filtro De Linha"
$ws.Cells.Item(12, 4).Value = 20

# Row 13
$ws.Cells.Item(13, 1).Value = 26
$ws.Cells.Item(13, 2).Value = "LIMPEZA"
$ws.Cells.Item(13, 3).Value = "Sabão Líquido"
$ws.Cells.Item(13, 4).Value = 800
Set-TextLiteral 13 5 "31/12/2025"

# Row 14
$ws.Cells.Item(14, 1).Value = 26
$ws.Cells.Item(14, 2).Value = "LIMPEZA"
$ws.Cells.Item(14, 3).Value = "Papel A1 Sulfite"
$ws.Cells.Item(14, 4).Value = 29
$ws.Cells.Item(14, 5).Value = ""

# Row 15
$ws.Cells.Item(15, 1).Value = 24
$ws.Cells.Item(15, 2).Value = "INFORMATICA"
$ws.Cells.Item(15, 3).Value = "Teclado Usb"
$ws.Cells.Item(15, 4).Value = 10
$ws.Cells.Item(15, 5).Value = ""

# Row 16
$ws.Cells.Item(16, 3).Value = "Caneta Vermelha"
$ws.Cells.Item(16, 4).Value = 99

# Row 17
$ws.Cells.Item(17, 1).Value = 29
$ws.Cells.Item(17, 2).Value = "DIVERSOS"
$ws.Cells.Item(17, 3).Value = "Café"
$ws.Cells.Item(17, 4).Value = 20
Set-TextLiteral 17 5 "19/06/2025"

# Row 18
$ws.Cells.Item(18, 1).Value = 29
$ws.Cells.Item(18, 2).Value = "DIVERSOS"
$ws.Cells.Item(18, 3).Value = "Açúcar"
$ws.Cells.Item(18, 4).Value = 49
Set-TextLiteral 18 5 "13/06/2025"

# Row 19
$ws.Cells.Item(19, 1).Value = 29
$ws.Cells.Item(19, 2).Value = "DIVERSOS"
$ws.Cells.Item(19, 3).Value = "Garrafa De Água"
$ws.Cells.Item(19, 4).Value = 19
Set-TextLiteral 19 5 "11/06/2025"

# Row 20
$ws.Cells.Item(20, 1).Value = 20
$ws.Cells.Item(20, 2).Value = "PAPELARIA"
$ws.Cells.Item(20, 3).Value = "Caneta Esferográfica Azul"
$ws.Cells.Item(20, 4).Value = 490
Set-TextLiteral 20 5 "31/12/2025"
